$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 11.915632
$ws.Cells.Item(2, 8).Value = 35.746896
$ws.Cells.Item(2, 9).Value = 0.2203762099850903
$ws.Cells.Item(2, 10).Value = 0.2203762099850904
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 103.4766596666667
$ws.Cells.Item(2, 14).Value = 310.429979
$ws.Cells.Item(2, 15).Value = 0.877785331764719
$ws.Cells.Item(2, 16).Value = 0.8777853317647188
$ws.Cells.Item(2, 17).Value = 1232.989797177243
$ws.Cells.Item(2, 18).Value = 11096.90817459518
$ws.Cells.Item(2, 19).Value = 0.1934430045948139
$ws.Cells.Item(2, 20).Value = 0.1934430045948139

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 11.915632
$ws.Cells.Item(3, 8).Value = 35.746896
$ws.Cells.Item(3, 9).Value = 0.2203762099850903
$ws.Cells.Item(3, 10).Value = 0.2203762099850904
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.8265796666666668
$ws.Cells.Item(3, 14).Value = 2.479739
$ws.Cells.Item(3, 15).Value = 0.007011818020336602
$ws.Cells.Item(3, 16).Value = 0.0070118180203366
$ws.Cells.Item(3, 17).Value = 9.849219126682669
$ws.Cells.Item(3, 18).Value = 88.64297214014401
$ws.Cells.Item(3, 19).Value = 0.001545237880426939
$ws.Cells.Item(3, 20).Value = 0.001545237880426939

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 11.915632
$ws.Cells.Item(4, 8).Value = 35.746896
$ws.Cells.Item(4, 9).Value = 0.2203762099850903
$ws.Cells.Item(4, 10).Value = 0.2203762099850904
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 13.58054833333333
$ws.Cells.Item(4, 14).Value = 40.741645
$ws.Cells.Item(4, 15).Value = 0.1152028502149446
$ws.Cells.Item(4, 16).Value = 0.1152028502149446
$ws.Cells.Item(4, 17).Value = 161.8208162982133
$ws.Cells.Item(4, 18).Value = 1456.38734668392
$ws.Cells.Item(4, 19).Value = 0.02538796750984954
$ws.Cells.Item(4, 20).Value = 0.02538796750984954

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 29.800487
$ws.Cells.Item(5, 8).Value = 89.40146100000001
$ws.Cells.Item(5, 9).Value = 0.5511514941691683
$ws.Cells.Item(5, 10).Value = 0.5511514941691684
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 103.4766596666667
$ws.Cells.Item(5, 14).Value = 310.429979
$ws.Cells.Item(5, 15).Value = 0.877785331764719
$ws.Cells.Item(5, 16).Value = 0.8777853317647188
$ws.Cells.Item(5, 17).Value = 3083.654851199924
$ws.Cells.Item(5, 18).Value = 27752.89366079932
$ws.Cells.Item(5, 19).Value = 0.483792697161904
$ws.Cells.Item(5, 20).Value = 0.483792697161904

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 29.800487
$ws.Cells.Item(6, 8).Value = 89.40146100000001
$ws.Cells.Item(6, 9).Value = 0.5511514941691683
$ws.Cells.Item(6, 10).Value = 0.5511514941691684
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.8265796666666668
$ws.Cells.Item(6, 14).Value = 2.479739
$ws.Cells.Item(6, 15).Value = 0.007011818020336602
$ws.Cells.Item(6, 16).Value = 0.0070118180203366
$ws.Cells.Item(6, 17).Value = 24.63247661096434
$ws.Cells.Item(6, 18).Value = 221.6922894986791
$ws.Cells.Item(6, 19).Value = 0.003864573978750818
$ws.Cells.Item(6, 20).Value = 0.003864573978750818

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 29.800487
$ws.Cells.Item(7, 8).Value = 89.40146100000001
$ws.Cells.Item(7, 9).Value = 0.5511514941691683
$ws.Cells.Item(7, 10).Value = 0.5511514941691684
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 13.58054833333333
$ws.Cells.Item(7, 14).Value = 40.741645
$ws.Cells.Item(7, 15).Value = 0.1152028502149446
$ws.Cells.Item(7, 16).Value = 0.1152028502149446
$ws.Cells.Item(7, 17).Value = 404.7069540603717
$ws.Cells.Item(7, 18).Value = 3642.362586543346
$ws.Cells.Item(7, 19).Value = 0.06349422302851361
$ws.Cells.Item(7, 20).Value = 0.06349422302851361

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 12.35338333333333
$ws.Cells.Item(8, 8).Value = 37.06015
$ws.Cells.Item(8, 9).Value = 0.2284722958457413
$ws.Cells.Item(8, 10).Value = 0.2284722958457413
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 103.4766596666667
$ws.Cells.Item(8, 14).Value = 310.429979
$ws.Cells.Item(8, 15).Value = 0.877785331764719
$ws.Cells.Item(8, 16).Value = 0.8777853317647188
$ws.Cells.Item(8, 17).Value = 1278.286842915206
$ws.Cells.Item(8, 18).Value = 11504.58158623685
$ws.Cells.Item(8, 19).Value = 0.200549630008001
$ws.Cells.Item(8, 20).Value = 0.200549630008001

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 12.35338333333333
$ws.Cells.Item(9, 8).Value = 37.06015
$ws.Cells.Item(9, 9).Value = 0.2284722958457413
$ws.Cells.Item(9, 10).Value = 0.2284722958457413
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.8265796666666668
$ws.Cells.Item(9, 14).Value = 2.479739
$ws.Cells.Item(9, 15).Value = 0.007011818020336602
$ws.Cells.Item(9, 16).Value = 0.0070118180203366
$ws.Cells.Item(9, 17).Value = 10.21105547787222
$ws.Cells.Item(9, 18).Value = 91.89949930085001
$ws.Cells.Item(9, 19).Value = 0.001602006161158844
$ws.Cells.Item(9, 20).Value = 0.001602006161158844

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 12.35338333333333
$ws.Cells.Item(10, 8).Value = 37.06015
$ws.Cells.Item(10, 9).Value = 0.2284722958457413
$ws.Cells.Item(10, 10).Value = 0.2284722958457413
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 13.58054833333333
$ws.Cells.Item(10, 14).Value = 40.741645
$ws.Cells.Item(10, 15).Value = 0.1152028502149446
$ws.Cells.Item(10, 16).Value = 0.1152028502149446
$ws.Cells.Item(10, 17).Value = 167.7657194385278
$ws.Cells.Item(10, 18).Value = 1509.89147494675
$ws.Cells.Item(10, 19).Value = 0.02632065967658144
$ws.Cells.Item(10, 20).Value = 0.02632065967658144
